# LOB1206.xlsx edit: remove the standalone "Docentes responsáveis" value row
# (row 13, which only carried B/C with no A label) and patch the handful of
# B/C cells whose text content changed as part of the same content refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete row 13 entirely (shifts rows 14-22 up to 13-21, and the sheet
#    dimension / used range shrinks from A1:C22 to A1:C21 automatically).
$ws.Rows.Item(13).Delete()

# 2) After the shift, a handful of B/C cells need their text corrected to
#    match the refreshed content.
$ws.Range("B10").Value = "5840942 - Marco Aurélio Kondracki de Alcântara"
$ws.Range("C10").Value = "5840942 - Marco Aurélio Kondracki de Alcântara"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B15").Value = "01/01/2020"
$ws.Range("C15").Value = "01/01/2020"

$ws.Range("B18").Value = "5840942 - Marco Aurélio Kondracki de Alcântara"
$ws.Range("C18").Value = "5840942 - Marco Aurélio Kondracki de Alcântara"

$ws.Range("B19").Value = "A avaliação será feita mediante duas avaliações escritas de igual peso (P1 e P2). Alternativamente, essas avaliações escritas poderão ser substituídas por duas notas de igual peso (NOTA 1 e NOTA 2). Essas NOTAS 1 e 2 serão dadas pela média entre atividades desenvolvidas em aula, trabalhos e relatórios de aulas práticas."
$ws.Range("C19").Value = "A avaliação será feita mediante duas avaliações escritas de igual peso (P1 e P2). Alternativamente, essas avaliações escritas poderão ser substituídas por duas notas de igual peso (NOTA 1 e NOTA 2). Essas NOTAS 1 e 2 serão dadas pela média entre atividades desenvolvidas em aula, trabalhos e relatórios de aulas práticas."

$ws.Range("B20").Value = "O aluno poderá optar por dois critérios de avaliação:Critério 1: NF = (P1+P2)/2; ouCritério 2: NF = (NOTA 1 + NOTA 2)/2Sendo P1 e P2 avaliações escritas e NOTA 1 e NOTA 2 obtidas em atividades desenvolvidas em aula, trabalhos e relatórios de aulas práticas."
$ws.Range("C20").Value = "O aluno poderá optar por dois critérios de avaliação:Critério 1: NF = (P1+P2)/2; ouCritério 2: NF = (NOTA 1 + NOTA 2)/2Sendo P1 e P2 avaliações escritas e NOTA 1 e NOTA 2 obtidas em atividades desenvolvidas em aula, trabalhos e relatórios de aulas práticas."

$ws.Range("B21").Value = "Exame Final (EF) para alunos com Nota Final (NF) maior ou igual a 3,0 e menor do que 6,5 e frequência superior a 70%. Será considerado aprovado o aluno que tenha obtido Média Final (MF) igual ou maior do que 5,0, sendo MF = (NF+EF)/2."
$ws.Range("C21").Value = "Exame Final (EF) para alunos com Nota Final (NF) maior ou igual a 3,0 e menor do que 6,5 e frequência superior a 70%. Será considerado aprovado o aluno que tenha obtido Média Final (MF) igual ou maior do que 5,0, sendo MF = (NF+EF)/2."
